$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire row for collaborator "FABIO FIGUEIRO MACHIDA" (row 12)
$ws.Rows.Item(12).Delete()

# Update the selection to match the post-edit state
$ws.Range("A12").Select()
